$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("publications")

$row = 18

$ws.Cells.Item($row, 1).Value = 44932
$ws.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item($row, 2).Value = "Nematic Torques in Scalar Active Matter: when Fluctuations Favor Polar Order and Persistence"
$ws.Cells.Item($row, 3).Value = "G. Spera, C. Duclut, M. Durand, J. Tailleur"
$ws.Cells.Item($row, 4).Value = "preprint"
$ws.Cells.Item($row, 5).Value = "/"
$ws.Cells.Item($row, 6).Value = "/"
$ws.Cells.Item($row, 7).Value = "We study the impact of nematic alignment on scalar active matter in the disordered phase. We show that nematic torques control the emergent physics of particles interacting via pairwise forces and can either induce or prevent phase separation. The underlying mechanism is a fluctuation-induced renormalization of the mass of the polar field that generically arises from nematic torques. The correlations between the fluctuations of the polar and nematic fields indeed conspire to increase the particle persistence length, contrary to what phenomenological computations predict. This effect is generic and our theory also quantitatively accounts for how nematic torques enhance particle accumulation along confining boundaries and opposes demixing in mixtures of active and passive particles. "
$ws.Cells.Item($row, 8).Value = "spera2023nematic"
$ws.Cells.Item($row, 9).Value = "/"

$arxivCell = $ws.Cells.Item($row, 10)
$arxivCell.NumberFormat = "@"
$arxivCell.Value = "2301.02568"
$arxivCell.Style = "Normal"

$ws.Cells.Item($row, 11).Value = "activeMatter, nematic"

$ws.Range("C27").Select()
